# OW-976, updated test results
# Appends 14 new rows (30-43) to the "Results" sheet, repeating the
# existing A/B/C (Rdata/timeLimit/callNumber) pattern from rows 16-29
# and filling in new runTime (D) values captured from the latest run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each new row's A/B/C values mirror the analogous row 14 rows above it
# (30<-16, 31<-17, ... 43<-29); only column D (runTime) is new data.
$newRows = @(
    @{ Row = 30; Template = 16; D = "0.0699999999999363" },
    @{ Row = 31; Template = 17; D = "0.0899999999999181" },
    @{ Row = 32; Template = 18; D = "0.0299999999999727" },
    @{ Row = 33; Template = 19; D = "0.170000000000073" },
    @{ Row = 34; Template = 20; D = "0.399999999999864" },
    @{ Row = 35; Template = 21; D = "0.210000000000036" },
    @{ Row = 36; Template = 22; D = "0.0499999999999545" },
    @{ Row = 37; Template = 23; D = "0.0599999999999454" },
    @{ Row = 38; Template = 24; D = "0.309999999999945" },
    @{ Row = 39; Template = 25; D = "0.0599999999999454" },
    @{ Row = 40; Template = 26; D = "0.0599999999999454" },
    @{ Row = 41; Template = 27; D = "0.160000000000082" },
    @{ Row = 42; Template = 28; D = "0.210000000000036" },
    @{ Row = 43; Template = 29; D = "0.200000000000045" }
)

foreach ($nr in $newRows) {
    $row = $nr.Row
    $tpl = $nr.Template

    $colA = $ws.Cells.Item($tpl, 1).Value2
    $colB = $ws.Cells.Item($tpl, 2).Value2
    $colC = $ws.Cells.Item($tpl, 3).Value2

    # Write values. A is already non-numeric text; B/C/D look like numbers
    # so they need a leading apostrophe to force text storage (matching the
    # shared-string cell type used throughout this table).
    $ws.Cells.Item($row, 1).Value = "'" + $colA
    $ws.Cells.Item($row, 2).Value = "'" + $colB
    $ws.Cells.Item($row, 3).Value = "'" + $colC
    $ws.Cells.Item($row, 4).Value = "'" + $nr.D

    # Re-apply the clean formatting (wrap text, General number format) from
    # the template row so the quote-prefix marker above doesn't linger.
    $ws.Cells.Item($tpl, 1).Copy()
    $ws.Cells.Item($row, 1).PasteSpecial(-4122)
    $ws.Cells.Item($tpl, 2).Copy()
    $ws.Cells.Item($row, 2).PasteSpecial(-4122)
    $ws.Cells.Item($tpl, 3).Copy()
    $ws.Cells.Item($row, 3).PasteSpecial(-4122)
    $ws.Cells.Item($tpl, 4).Copy()
    $ws.Cells.Item($row, 4).PasteSpecial(-4122)
}

$ws.Range("A1").Select()
